# Update "想去人数" (F column) values across the four sheets to reflect
# newly scraped counts (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 261
$ws1.Range("F5").Value  = 1037
$ws1.Range("F6").Value  = 2302
$ws1.Range("F7").Value  = 224
$ws1.Range("F8").Value  = 668
$ws1.Range("F9").Value  = 33
$ws1.Range("F10").Value = 204
$ws1.Range("F11").Value = 163
$ws1.Range("F12").Value = 664
$ws1.Range("F13").Value = 60
$ws1.Range("F15").Value = 1356
$ws1.Range("F18").Value = 190

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value  = 22
$ws2.Range("F6").Value  = 13
$ws2.Range("F11").Value = 37
$ws2.Range("F12").Value = 212

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6332
$ws3.Range("F5").Value = 220

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 6332
$ws4.Range("F6").Value  = 220
$ws4.Range("F10").Value = 22
$ws4.Range("F12").Value = 261
$ws4.Range("F13").Value = 1037
$ws4.Range("F14").Value = 13
$ws4.Range("F17").Value = 2302
$ws4.Range("F19").Value = 224
$ws4.Range("F21").Value = 37
$ws4.Range("F22").Value = 668
$ws4.Range("F23").Value = 33
$ws4.Range("F24").Value = 204
$ws4.Range("F25").Value = 212
$ws4.Range("F26").Value = 163
$ws4.Range("F27").Value = 664
$ws4.Range("F28").Value = 60
$ws4.Range("F31").Value = 1356
$ws4.Range("F36").Value = 190
